$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "oficina de Bankia receptora" -> "oficina de Grupo CaixaBank receptora"
Replace-Text "La Reserva deberá formalizarse en la oficina de Bankia receptora de la presente comunicación, en el plazo máximo de 7 días hábiles a contar desde la fecha de recepción del documento de Reserva." "La Reserva deberá formalizarse en la oficina de Grupo CaixaBank receptora de la presente comunicación, en el plazo máximo de 7 días hábiles a contar desde la fecha de recepción del documento de Reserva."

# 2. merge runs around "documento de reserva ... Oficina Bankia"
Replace-Text "Es imprescindible que el documento de reserva sea firmado por duplicado por todos y cada uno de los compradores que deberán coincidir con los titulares de la propuesta aprobada, así como por parte de los apoderados de la Oficina Bankia. Cualquier modificación en el documento deberá autorizarse expresamente por el Área de " "Es imprescindible que el documento de reserva sea firmado por duplicado por todos y cada uno de los compradores que deberán coincidir con los titulares de la propuesta aprobada, así como por parte de los apoderados de la Oficina Grupo CaixaBank. Cualquier modificación en el documento deberá autorizarse expresamente por el Área de "

# 3. merge runs (Reiteramos ... documento de reserva, implicará ...)
Replace-Text "Reiteramos que los plazos estipulados en el presente documento se consideran esenciales para la toma de decisión indicada por lo que tanto; (i) la ausencia de comunicación del Interesado para formalizar el documento de Reserva, como; (ii) la ausencia de formalización del documento de reserva, implicará la paralización en los trámites de la operación de venta, estando en disposición de la Propietaria del Inmueble de anular la presente oferta, y proceder a la tramitación de la siguiente. " "Reiteramos que los plazos estipulados en el presente documento se consideran esenciales para la toma de decisión indicada por lo que tanto; (i) la ausencia de comunicación del Interesado para formalizar el documento de Reserva, como; (ii) la ausencia de formalización del documento de reserva, implicará la paralización en los trámites de la operación de venta, estando en disposición de la Propietaria del Inmueble de anular la presente oferta, y proceder a la tramitación de la siguiente. "

# 4. merge + shorten "ALTA EN NEO CLIENTES" paragraph, drop "Habitat" sentence
Replace-Text "La Oficina, deberá dar de alta como cliente Bankia al comprador/es de la operación si la oferta aprobada corresponde un activo/lote cuyo propietario sea Bankia, o bien, en caso de corresponder a propietario distinto de Bankia, darlo/s de alta como cliente Bankia Habitat." "La Oficina, deberá dar de alta como cliente Grupo CaixaBank al comprador/es de la operación si la oferta aprobada corresponde un activo/lote cuyo propietario sea Grupo CaixaBank."

# 5. "cuenta en Bankia" -> "cuenta en Grupo CaixaBank"
Replace-Text "Si el comprador no dispone de cuenta en Bankia y/o financia otra entidad, el cobro se hará mediante cheque bancario." "Si el comprador no dispone de cuenta en Grupo CaixaBank y/o financia otra entidad, el cobro se hará mediante cheque bancario."

# 6. "financiada por Bankia o es cliente de Bankia" -> "...Grupo CaixaBank..."
Replace-Text "Si se trata de una operación financiada por Bankia o es cliente de Bankia, se realizará " "Si se trata de una operación financiada por Grupo CaixaBank o es cliente de Grupo CaixaBank, se realizará "

# 7. merge "Asimismo," + " les indicamos..." runs
Replace-Text "Asimismo, les indicamos que la operación de venta ha sido aprobada por el Comité de la Entidad en los términos recogidos en la presente comunicación en cuanto a titulares y demás condiciones. " "Asimismo, les indicamos que la operación de venta ha sido aprobada por el Comité de la Entidad en los términos recogidos en la presente comunicación en cuanto a titulares y demás condiciones. "

# 8. merge "Quedamos ... cordiales." + " " runs
Replace-Text "Quedamos a su disposición para cualquier consulta o aclaración. Saludos cordiales. " "Quedamos a su disposición para cualquier consulta o aclaración. Saludos cordiales. "
